$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price / 1h-volume snapshot.
# Two pairs of rows also swapped rank position: Toncoin/Monero (25-26)
# and TrustWalletToken/Hedera (38-39), so their Coin/Link/Price/Volume
# cells are rewritten together.
# Price cells (column D) are plain-text in the workbook (e.g. "27.102.48"
# uses "." as a thousands separator), so every Price write is prefixed with
# a leading apostrophe to force text entry -- otherwise Excel would parse
# strings such as "1.009" or "92.60" as numbers and mangle/round them.

$ws.Range("D2").Value = "'27.102.48"
$ws.Range("E2").Value = '  -1.14%  '
$ws.Range("D3").Value = "'1.822.93"
$ws.Range("E3").Value = '  -1.18%  '
$ws.Range("D4").Value = "'1.009"
$ws.Range("E4").Value = '  -0.56%  '
$ws.Range("D5").Value = "'311.69"
$ws.Range("E5").Value = '  -1.51%  '
$ws.Range("D6").Value = "'1.009"
$ws.Range("E6").Value = '  -0.50%  '
$ws.Range("D7").Value = "'0.4631"
$ws.Range("E7").Value = '  -2.37%  '
$ws.Range("D8").Value = "'0.3639"
$ws.Range("E8").Value = '  -1.76%  '
$ws.Range("D9").Value = "'0.07296"
$ws.Range("E9").Value = '  -2.23%  '
$ws.Range("D10").Value = "'0.8702"
$ws.Range("E10").Value = '  -2.03%  '
$ws.Range("D11").Value = "'20.13"
$ws.Range("E11").Value = '  -1.95%  '
$ws.Range("D12").Value = "'1.872.65"
$ws.Range("E12").Value = '  +0.76%  '
$ws.Range("D13").Value = "'0.07583"
$ws.Range("E13").Value = '  +2.41%  '
$ws.Range("D14").Value = "'5.346"
$ws.Range("E14").Value = '  -2.72%  '
$ws.Range("D15").Value = "'92.60"
$ws.Range("E15").Value = '  -0.84%  '
$ws.Range("D16").Value = "'6.475"
$ws.Range("E16").Value = '  -1.86%  '
$ws.Range("D17").Value = "'1.008"
$ws.Range("E17").Value = '  -0.67%  '
$ws.Range("D18").Value = "'0.000008645"
$ws.Range("E18").Value = '  -2.53%  '
$ws.Range("E19").Value = '  -0.56%  '
$ws.Range("D20").Value = "'27.394.57"
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("D21").Value = "'14.47"
$ws.Range("D22").Value = "'5.198"
$ws.Range("E22").Value = '  -2.75%  '
$ws.Range("D23").Value = "'10.56"
$ws.Range("E23").Value = '  -1.66%  '
$ws.Range("D24").Value = "'2.093.75"
$ws.Range("E24").Value = '  +0.87%  '
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = "'1.876"
$ws.Range("E25").Value = '  -1.88%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = "'151.71"
$ws.Range("E26").Value = '  -0.55%  '
$ws.Range("D27").Value = "'18.26"
$ws.Range("E27").Value = '  -2.39%  '
$ws.Range("D28").Value = "'2.100"
$ws.Range("E28").Value = '  -3.77%  '
$ws.Range("D29").Value = "'116.23"
$ws.Range("E29").Value = '  -1.60%  '
$ws.Range("D30").Value = "'5.059"
$ws.Range("E30").Value = '  -4.37%  '
$ws.Range("D31").Value = "'0.08915"
$ws.Range("E31").Value = '  -0.72%  '
$ws.Range("D32").Value = "'2.963"
$ws.Range("E32").Value = '  +0.35%  '
$ws.Range("D33").Value = "'0.7346"
$ws.Range("E33").Value = '  -3.79%  '
$ws.Range("D34").Value = "'4.453"
$ws.Range("E34").Value = '  -2.63%  '
$ws.Range("D35").Value = "'1.139"
$ws.Range("D36").Value = "'1.010"
$ws.Range("E36").Value = '  -0.37%  '
$ws.Range("D37").Value = "'2.519"
$ws.Range("E37").Value = '  +4.96%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = "'1.073"
$ws.Range("E38").Value = '  -3.13%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = "'0.05258"
$ws.Range("E39").Value = '  -2.03%  '
$ws.Range("D40").Value = "'0.01922"
$ws.Range("E40").Value = '  -2.43%  '
$ws.Range("D41").Value = "'2.936"
$ws.Range("E41").Value = '  -2.42%  '
$ws.Range("D42").Value = "'7.158"
$ws.Range("E42").Value = '  -2.43%  '
$ws.Range("D43").Value = "'0.5211"
$ws.Range("E43").Value = '  -2.95%  '
$ws.Range("D44").Value = "'0.1632"
$ws.Range("E44").Value = '  -2.18%  '
$ws.Range("D45").Value = "'8.273"
$ws.Range("E45").Value = '  -3.34%  '
$ws.Range("D46").Value = "'0.4888"
$ws.Range("E46").Value = '  -1.58%  '
$ws.Range("D47").Value = "'1.009"
$ws.Range("E47").Value = '  -0.58%  '
$ws.Range("D48").Value = "'10.17"
$ws.Range("E48").Value = '  -3.79%  '
$ws.Range("D49").Value = "'103.93"
$ws.Range("E49").Value = '  -1.17%  '
$ws.Range("D50").Value = "'1.634"
$ws.Range("E50").Value = '  -3.04%  '
$ws.Range("D51").Value = "'0.06255"
$ws.Range("E51").Value = '  -1.32%  '
